$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for "Haba" at Vega Central
# Mapocho de Santiago. Insert a fresh row at row 99 (pushing the existing
# rows 99-139 down to 100-140) and populate it with the new record.
$ws.Rows.Item(99).Insert()

$ws.Range("A99").Value = 9
$ws.Range("B99").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C99").Value = "Metropolitana"
$ws.Range("D99").Value = 44466
$ws.Range("E99").Value = 13
$ws.Range("F99").Value = 100112026
$ws.Range("G99").Value = "Haba"
$ws.Range("H99").Value = "Sin especificar"
$ws.Range("I99").Value = "Primera"
$ws.Range("J99").Value = 43
$ws.Range("K99").Value = 14000
$ws.Range("L99").Value = 15000
$ws.Range("M99").Value = 14512
$ws.Range("N99").Value = '$/saco 25 kilos'
$ws.Range("O99").Value = "Provincia de Limarí"
$ws.Range("P99").Value = 580
$ws.Range("Q99").Value = 25
$ws.Range("R99").Value = "Hortaliza"
